# New polymers added to the "Structures" sheet: PM6, P5TCN-r, P5TCN
# (rows 99-101), each with its Name plus dD / dP / dH (Hansen solubility
# parameter) values. Columns B (Image), C (SMILES) and D (BigSMILES) are
# intentionally left blank for these entries, same as in the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- names first --------------------------------------------------------
$ws.Range("A99").Value  = "PM6"
$ws.Range("A100").Value = "P5TCN-r"
$ws.Range("A101").Value = "P5TCN"

# Row 99 should look like every other Name cell in column A: centered,
# wrapped text. Copy that look straight from the cell above instead of
# re-deriving it property-by-property.
$ws.Range("A98").Copy()
$ws.Range("A99").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Rows 100/101 use a new look: centered both ways, but NOT wrapped. Build
# it once on A100 ...
$ws.Range("A100").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A100").VerticalAlignment = -4108     # xlCenter

# ... then reuse the exact same style for A101 via a format copy.
$ws.Range("A100").Copy()
$ws.Range("A101").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ---- row heights (matches the rest of the sheet) ------------------------
$ws.Rows.Item(99).RowHeight = 90
$ws.Rows.Item(100).RowHeight = 90
$ws.Rows.Item(101).RowHeight = 90

# ---- dD / dP / dH values -------------------------------------------------
$ws.Range("E99").Value = 17.1
$ws.Range("F99").Value = 1.7
$ws.Range("G99").Value = 4.29

$ws.Range("E100").Value = 16.9
$ws.Range("F100").Value = 3.4
$ws.Range("G100").Value = 5.5

$ws.Range("E101").Value = 17.7
$ws.Range("F101").Value = 3.75
$ws.Range("G101").Value = 4.85

# ---- leave the selection where the author left it ------------------------
$ws.Range("C99").Select() | Out-Null
